$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (labels are unchanged; shared-string indices will be
#     renumbered once the old demo rows below are replaced) ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "use_sublocation"
$ws.Range("C1").Value = "sublocation_min"
$ws.Range("D1").Value = "sublocation_max"

# --- New demo data (stock_location.name made consistent with
#     pv_pn.PNUser9), entered in the order it was originally keyed in,
#     i.e. before the alphabetical sort performed below ---
$names  = @("FINGOODS","ELEC-PCB","ELEC-COMPS","ELEC-WIRE","DOCS","SHOPSUP","MECH","SHIPPING","ELEC-CONNS","MECH-HW","YARD")
$useSub = @(1,1,1,1,1,1,1,1,1,1,0)
$subMin = @(1,1,1,1,1,1,1,1,1,1,0)
$subMax = @(10,10,20,10,10,10,10,10,10,10,0)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = $useSub[$i]
    $ws.Cells.Item($r, 3).Value = $subMin[$i]
    $ws.Cells.Item($r, 4).Value = $subMax[$i]
}

# --- Sort the bulk of the list (rows 2-11) alphabetically by name;
#     the catch-all "YARD" row stays pinned at the bottom, outside the
#     sorted range, just like "stock room 2" was before ---
$ws.Range("A2:D11").Sort($ws.Range("A2"), 1)

# --- Re-apply the centred number/flag formatting to the whole data body
#     (column A keeps the default/general style) ---
$ws.Range("B2:D12").HorizontalAlignment = -4108

# --- Record the sort state Excel keeps on the worksheet after an
#     in-place column sort ---
$sortFields = $ws.Sort.SortFields
$sortFields.Clear()
$sortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:A11"))
$ws.Sort.Header = 2
$ws.Sort.Apply()
